# Correction to Asn flux: remove the "Asparagine pos" column (column F)
# entirely from the worksheet. This shifts every column to its right
# (Asparagine-13C4 pos .. Valine pos, and all numeric data beneath them)
# one place to the left, and drops the now-unused "Asparagine pos" shared
# string from the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Columns.Item(6).Delete()
